# This script applies a data re-sort among several groups of rows in the
# "artfynd" worksheet. Within each group, the observation metadata
# (species id/name, taxon id, Swedish name, scientific name, author,
# coordinates, and - where relevant - age/stage/activity/method/comment
# columns) is rotated/swapped between the rows belonging to that group,
# while the columns that are common to the whole group (locality, date,
# municipality, reporter, etc.) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowValues($sheet, $row, $cols) {
    $vals = @{}
    foreach ($col in $cols) {
        $vals[$col] = $sheet.Range("$col$row").Value2
    }
    return $vals
}

function Set-RowValues($sheet, $row, $cols, $vals) {
    foreach ($col in $cols) {
        $v = $vals[$col]
        if ($v -eq $null) { $v = "" }
        $sheet.Range("$col$row").Value2 = $v
    }
}

# --- Standard column set used by most groups ---
$cols = @("A","B","E","F","G","H","Q","R")

# --- Group 1: rows 2, 3, 5 (cyclic shift: 2 -> 3 -> 5 -> 2) ---
$v2 = Get-RowValues $ws 2 $cols
$v3 = Get-RowValues $ws 3 $cols
$v5 = Get-RowValues $ws 5 $cols
Set-RowValues $ws 3 $cols $v2
Set-RowValues $ws 5 $cols $v3
Set-RowValues $ws 2 $cols $v5

# --- Group 2: rows 6, 7, 10 (cyclic shift: 6 -> 7 -> 10 -> 6) ---
$v6 = Get-RowValues $ws 6 $cols
$v7 = Get-RowValues $ws 7 $cols
$v10 = Get-RowValues $ws 10 $cols
Set-RowValues $ws 7 $cols $v6
Set-RowValues $ws 10 $cols $v7
Set-RowValues $ws 6 $cols $v10

# --- Group 3: rows 17, 18 (simple swap) ---
$v17 = Get-RowValues $ws 17 $cols
$v18 = Get-RowValues $ws 18 $cols
Set-RowValues $ws 17 $cols $v18
Set-RowValues $ws 18 $cols $v17

# --- Group 4: rows 34, 36 (simple swap) ---
$v34 = Get-RowValues $ws 34 $cols
$v36 = Get-RowValues $ws 36 $cols
Set-RowValues $ws 34 $cols $v36
Set-RowValues $ws 36 $cols $v34

# --- Group 5: rows 61, 62 (simple swap, extra columns K/L/M/N/AC) ---
$cols61 = @("A","B","E","F","G","H","K","L","M","N","Q","R","AC")
$v61 = Get-RowValues $ws 61 $cols61
$v62 = Get-RowValues $ws 62 $cols61
Set-RowValues $ws 61 $cols61 $v62
Set-RowValues $ws 62 $cols61 $v61

# --- Group 6: rows 69, 70 (simple swap, includes D column) ---
# Note: column H is special-cased here - per the source data, row 69's
# Auktor is overwritten with row 70's old Auktor value ("Tibell"), but
# row 70's Auktor cell is left as-is (it keeps reading "Tibell" too).
$cols69 = @("A","B","D","E","F","G","Q","R")
$v69 = Get-RowValues $ws 69 $cols69
$v70 = Get-RowValues $ws 70 $cols69
$h70 = $ws.Range("H70").Value2
Set-RowValues $ws 69 $cols69 $v70
Set-RowValues $ws 70 $cols69 $v69
$ws.Range("H69").Value2 = $h70
